$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"14.04383033333333"
$ws.Range("H2").Value = [double]"42.131491"
$ws.Range("I2").Value = [double]"0.158930310642385"
$ws.Range("J2").Value = [double]"0.158930310642385"
$ws.Range("M2").Value = [double]"28.22141"
$ws.Range("N2").Value = [double]"84.66423"
$ws.Range("O2").Value = [double]"0.007710741921554872"
$ws.Range("P2").Value = [double]"0.007710741921554872"
$ws.Range("Q2").Value = [double]"396.3366938074367"
$ws.Range("R2").Value = [double]"3567.03024426693"
$ws.Range("S2").Value = [double]"0.001225470608875976"
$ws.Range("T2").Value = [double]"0.001225470608875976"
$ws.Range("G3").Value = [double]"14.04383033333333"
$ws.Range("H3").Value = [double]"42.131491"
$ws.Range("I3").Value = [double]"0.158930310642385"
$ws.Range("J3").Value = [double]"0.158930310642385"
$ws.Range("O3").Value = [double]"0.001803104886918205"
$ws.Range("P3").Value = [double]"0.001803104886918206"
$ws.Range("Q3").Value = [double]"92.68065728817533"
$ws.Range("R3").Value = [double]"834.1259155935779"
$ws.Range("S3").Value = [double]"0.0002865680197987129"
$ws.Range("T3").Value = [double]"0.0002865680197987128"
$ws.Range("G4").Value = [double]"14.04383033333333"
$ws.Range("H4").Value = [double]"42.131491"
$ws.Range("I4").Value = [double]"0.158930310642385"
$ws.Range("J4").Value = [double]"0.158930310642385"
$ws.Range("M4").Value = [double]"1.757996666666666"
$ws.Range("N4").Value = [double]"5.27399"
$ws.Range("O4").Value = [double]"0.000480325348578274"
$ws.Range("P4").Value = [double]"0.0004803253485782741"
$ws.Range("Q4").Value = [double]"24.68900691323222"
$ws.Range("R4").Value = [double]"222.20106221909"
$ws.Range("S4").Value = [double]"7.633825685895694E-05"
$ws.Range("T4").Value = [double]"7.633825685895694E-05"
$ws.Range("G5").Value = [double]"14.04383033333333"
$ws.Range("H5").Value = [double]"42.131491"
$ws.Range("I5").Value = [double]"0.158930310642385"
$ws.Range("J5").Value = [double]"0.158930310642385"
$ws.Range("M5").Value = [double]"3623.433471666667"
$ws.Range("N5").Value = [double]"10870.300415"
$ws.Range("O5").Value = [double]"0.9900058278429487"
$ws.Range("P5").Value = [double]"0.9900058278429487"
$ws.Range("Q5").Value = [double]"50886.88490020764"
$ws.Range("R5").Value = [double]"457981.9641018687"
$ws.Range("S5").Value = [double]"0.1573419337568513"
$ws.Range("T5").Value = [double]"0.1573419337568513"
$ws.Range("I6").Value = [double]"0.3128977537755355"
$ws.Range("J6").Value = [double]"0.3128977537755354"
$ws.Range("M6").Value = [double]"28.22141"
$ws.Range("N6").Value = [double]"84.66423"
$ws.Range("O6").Value = [double]"0.007710741921554872"
$ws.Range("P6").Value = [double]"0.007710741921554872"
$ws.Range("Q6").Value = [double]"780.2971046235167"
$ws.Range("R6").Value = [double]"7022.673941611651"
$ws.Range("S6").Value = [double]"0.002412673827197375"
$ws.Range("T6").Value = [double]"0.002412673827197375"
$ws.Range("I7").Value = [double]"0.3128977537755355"
$ws.Range("J7").Value = [double]"0.3128977537755354"
$ws.Range("O7").Value = [double]"0.001803104886918205"
$ws.Range("P7").Value = [double]"0.001803104886918206"
$ws.Range("S7").Value = [double]"0.0005641874689383974"
$ws.Range("T7").Value = [double]"0.0005641874689383973"
$ws.Range("I8").Value = [double]"0.3128977537755355"
$ws.Range("J8").Value = [double]"0.3128977537755354"
$ws.Range("M8").Value = [double]"1.757996666666666"
$ws.Range("N8").Value = [double]"5.27399"
$ws.Range("O8").Value = [double]"0.000480325348578274"
$ws.Range("P8").Value = [double]"0.0004803253485782741"
$ws.Range("Q8").Value = [double]"48.60705786627221"
$ws.Range("R8").Value = [double]"437.46352079645"
$ws.Range("S8").Value = [double]"0.000150292722651593"
$ws.Range("T8").Value = [double]"0.000150292722651593"
$ws.Range("I9").Value = [double]"0.3128977537755355"
$ws.Range("J9").Value = [double]"0.3128977537755354"
$ws.Range("M9").Value = [double]"3623.433471666667"
$ws.Range("N9").Value = [double]"10870.300415"
$ws.Range("O9").Value = [double]"0.9900058278429487"
$ws.Range("P9").Value = [double]"0.9900058278429487"
$ws.Range("Q9").Value = [double]"100184.7408310725"
$ws.Range("R9").Value = [double]"901662.6674796523"
$ws.Range("S9").Value = [double]"0.3097705997567481"
$ws.Range("T9").Value = [double]"0.309770599756748"
$ws.Range("G10").Value = [double]"25.89747433333334"
$ws.Range("H10").Value = [double]"77.69242300000001"
$ws.Range("I10").Value = [double]"0.2930748622675039"
$ws.Range("J10").Value = [double]"0.2930748622675038"
$ws.Range("M10").Value = [double]"28.22141"
$ws.Range("N10").Value = [double]"84.66423"
$ws.Range("O10").Value = [double]"0.007710741921554872"
$ws.Range("P10").Value = [double]"0.007710741921554872"
$ws.Range("Q10").Value = [double]"730.8632411254767"
$ws.Range("R10").Value = [double]"6577.769170129291"
$ws.Range("S10").Value = [double]"0.002259824626639962"
$ws.Range("T10").Value = [double]"0.002259824626639962"
$ws.Range("G11").Value = [double]"25.89747433333334"
$ws.Range("H11").Value = [double]"77.69242300000001"
$ws.Range("I11").Value = [double]"0.2930748622675039"
$ws.Range("J11").Value = [double]"0.2930748622675038"
$ws.Range("O11").Value = [double]"0.001803104886918205"
$ws.Range("P11").Value = [double]"0.001803104886918206"
$ws.Range("Q11").Value = [double]"170.9074295507593"
$ws.Range("R11").Value = [double]"1538.166865956834"
$ws.Range("S11").Value = [double]"0.0005284447163874163"
$ws.Range("T11").Value = [double]"0.0005284447163874162"
$ws.Range("G12").Value = [double]"25.89747433333334"
$ws.Range("H12").Value = [double]"77.69242300000001"
$ws.Range("I12").Value = [double]"0.2930748622675039"
$ws.Range("J12").Value = [double]"0.2930748622675038"
$ws.Range("M12").Value = [double]"1.757996666666666"
$ws.Range("N12").Value = [double]"5.27399"
$ws.Range("O12").Value = [double]"0.000480325348578274"
$ws.Range("P12").Value = [double]"0.0004803253485782741"
$ws.Range("Q12").Value = [double]"45.52767355308555"
$ws.Range("R12").Value = [double]"409.74906197777"
$ws.Range("S12").Value = [double]"0.0001407712853781685"
$ws.Range("T12").Value = [double]"0.0001407712853781685"
$ws.Range("G13").Value = [double]"25.89747433333334"
$ws.Range("H13").Value = [double]"77.69242300000001"
$ws.Range("I13").Value = [double]"0.2930748622675039"
$ws.Range("J13").Value = [double]"0.2930748622675038"
$ws.Range("M13").Value = [double]"3623.433471666667"
$ws.Range("N13").Value = [double]"10870.300415"
$ws.Range("O13").Value = [double]"0.9900058278429487"
$ws.Range("P13").Value = [double]"0.9900058278429487"
$ws.Range("Q13").Value = [double]"93837.7753310284"
$ws.Range("R13").Value = [double]"844539.9779792556"
$ws.Range("S13").Value = [double]"0.2901458216390984"
$ws.Range("T13").Value = [double]"0.2901458216390983"
$ws.Range("G14").Value = [double]"20.774284"
$ws.Range("H14").Value = [double]"62.322852"
$ws.Range("I14").Value = [double]"0.2350970733145757"
$ws.Range("J14").Value = [double]"0.2350970733145757"
$ws.Range("M14").Value = [double]"28.22141"
$ws.Range("N14").Value = [double]"84.66423"
$ws.Range("O14").Value = [double]"0.007710741921554872"
$ws.Range("P14").Value = [double]"0.007710741921554872"
$ws.Range("Q14").Value = [double]"586.27958622044"
$ws.Range("R14").Value = [double]"5276.51627598396"
$ws.Range("S14").Value = [double]"0.001812772858841558"
$ws.Range("T14").Value = [double]"0.001812772858841558"
$ws.Range("G15").Value = [double]"20.774284"
$ws.Range("H15").Value = [double]"62.322852"
$ws.Range("I15").Value = [double]"0.2350970733145757"
$ws.Range("J15").Value = [double]"0.2350970733145757"
$ws.Range("O15").Value = [double]"0.001803104886918205"
$ws.Range("P15").Value = [double]"0.001803104886918206"
$ws.Range("Q15").Value = [double]"137.097518989624"
$ws.Range("R15").Value = [double]"1233.877670906616"
$ws.Range("S15").Value = [double]"0.000423904681793679"
$ws.Range("T15").Value = [double]"0.000423904681793679"
$ws.Range("G16").Value = [double]"20.774284"
$ws.Range("H16").Value = [double]"62.322852"
$ws.Range("I16").Value = [double]"0.2350970733145757"
$ws.Range("J16").Value = [double]"0.2350970733145757"
$ws.Range("M16").Value = [double]"1.757996666666666"
$ws.Range("N16").Value = [double]"5.27399"
$ws.Range("O16").Value = [double]"0.000480325348578274"
$ws.Range("P16").Value = [double]"0.0004803253485782741"
$ws.Range("Q16").Value = [double]"36.52112202438666"
$ws.Range("R16").Value = [double]"328.69009821948"
$ws.Range("S16").Value = [double]"0.0001129230836895556"
$ws.Range("T16").Value = [double]"0.0001129230836895556"
$ws.Range("G17").Value = [double]"20.774284"
$ws.Range("H17").Value = [double]"62.322852"
$ws.Range("I17").Value = [double]"0.2350970733145757"
$ws.Range("J17").Value = [double]"0.2350970733145757"
$ws.Range("M17").Value = [double]"3623.433471666667"
$ws.Range("N17").Value = [double]"10870.300415"
$ws.Range("O17").Value = [double]"0.9900058278429487"
$ws.Range("P17").Value = [double]"0.9900058278429487"
$ws.Range("Q17").Value = [double]"75274.23599550928"
$ws.Range("R17").Value = [double]"677468.1239595836"
$ws.Range("S17").Value = [double]"0.2327474726902509"
$ws.Range("T17").Value = [double]"0.2327474726902509"
